# ---------------------------------------------------------------------------
# Applies three content fixes to the T.L.E Reviewer document:
#   1. "Appriopriate" -> "Appropriate"  (also drops the stale spell-check
#      proofErr bookmarks that bracketed the old misspelling)
#   2. "Carbin" -> "Carbon", with the corrected word split across three runs
#      (Word's autocorrect/track-changes split: "Carb" / "o" / "n ...")
#   3. Drops the trailing "(Not Included in Exam)" qualifier from the
#      "Calculating the Cost of Production" heading.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParagraphContaining {
    param([string]$needle)
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Appriopriate -> Appropriate, and remove the now-stale proofErr markers
#    that wrapped the misspelling. proofErr bookmarks are anchored to the
#    paragraph itself, so we delete the whole paragraph (pruning them) and
#    rebuild it with the corrected text, re-splitting it back into the same
#    two runs ("Appropriate" / " for small kitchens.") it originally had by
#    briefly toggling a character property off/on across the boundary.
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphContaining "Appriopriate"
if ($p1 -ne $null) {
    $full = $p1.Range.Text
    $full = $full.Replace("Appriopriate", "Appropriate")
    # full includes the trailing paragraph mark already (Range.Text of a
    # paragraph always ends with \r)

    $pStart = $p1.Range.Start
    $pEnd = $p1.Range.End

    $d.Range($pStart, $pEnd).Delete()

    $insPt = $d.Range($pStart, $pStart)
    $insPt.Text = $full

    $splitAt = $pStart + 11    # length of "Appropriate"
    $rA = $d.Range($pStart, $splitAt)
    $rA.Font.Bold = 1
    $rA.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 2) Carbin -> Carbon, split as "Carb" / "o" / "n Monoxide ..." (all three
#    runs keep the existing bold formatting of the source run).
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Carbin", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Carbon", 2)

$p2 = Find-ParagraphContaining "Emission of Carbon"
if ($p2 -ne $null) {
    $pStart = $p2.Range.Start
    $b1 = $pStart + 16   # end of "Emission of Carb"
    $b2 = $pStart + 17   # end of "Emission of Carbo"

    $rA = $d.Range($pStart, $b1)
    $rB = $d.Range($b1, $b2)

    $rA.Font.Bold = 0
    $rA.Font.Bold = 1
    $rB.Font.Bold = 0
    $rB.Font.Bold = 1
}

# ---------------------------------------------------------------------------
# 3) Drop the "(Not Included in Exam)" qualifier.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute( `
    "Calculating the Cost of Production (Not Included in Exam)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Calculating the Cost of Production", 2)

Write-Output "done"
